# Update for 2024 November data release
# Adds new UI<->database mapping rows (MultiplexMicroscopy.tissue_fixative,
# MultiplexMicroscopy.imaging_assay_type, proteomic.* fields, image.organ_or_tissue)
# to the "Mapping" and "Must have properties" sheets, and logs the change on
# the "Change Log" sheet.

$wb = $excel.ActiveWorkbook

$wsMap   = $wb.Worksheets.Item("Mapping")
$wsMust  = $wb.Worksheets.Item("Must have properties")
$wsLog   = $wb.Worksheets.Item("Change Log")

# ---------------------------------------------------------------------------
# 1. "Mapping" sheet - six new data rows appended after the existing table
#    (rows 90-95). Columns: A=Page, B=Area, C=Display Name (blank), D=Full
#    Name, E=Node Name, F=Property Name.
# ---------------------------------------------------------------------------

$wsMap.Range("A90").Value = "DATA"
$wsMap.Range("B90").Value = "Widgets"
$wsMap.Range("D90").Value = "MultiplexMicroscopy.tissue_fixative"
$wsMap.Range("E90").Value = "MultiplexMicroscopy"
$wsMap.Range("F90").Value = "tissue_fixative"

$wsMap.Range("A91").Value = "DATA"
$wsMap.Range("B91").Value = "Widgets"
$wsMap.Range("D91").Value = "MultiplexMicroscopy.imaging_assay_type"
$wsMap.Range("E91").Value = "MultiplexMicroscopy"
$wsMap.Range("F91").Value = "imaging_assay_type"

$wsMap.Range("A92").Value = "DATA"
$wsMap.Range("B92").Value = "Widgets"
$wsMap.Range("D92").Value = "proteomic.analytical_fractions"
$wsMap.Range("E92").Value = "proteomic"
$wsMap.Range("F92").Value = "analytical_fractions"

$wsMap.Range("A93").Value = "DATA"
$wsMap.Range("B93").Value = "Widgets"
$wsMap.Range("D93").Value = "proteomic.instrument_make"
$wsMap.Range("E93").Value = "proteomic"
$wsMap.Range("F93").Value = "instrument_make"

$wsMap.Range("A94").Value = "DATA"
$wsMap.Range("B94").Value = "Widgets"
$wsMap.Range("D94").Value = "proteomic.proteomic_design_description"
$wsMap.Range("E94").Value = "proteomic"
$wsMap.Range("F94").Value = "proteomic_design_description"

$wsMap.Range("A95").Value = "DATA"
$wsMap.Range("B95").Value = "Widgets"
$wsMap.Range("D95").Value = "image.organ_or_tissue"
$wsMap.Range("E95").Value = "image"
$wsMap.Range("F95").Value = "organ_or_tissue"

$wsMap.Rows.Item(90).RowHeight = 19
$wsMap.Rows.Item(91).RowHeight = 19
$wsMap.Rows.Item(92).RowHeight = 18
$wsMap.Rows.Item(93).RowHeight = 18
$wsMap.Rows.Item(94).RowHeight = 18
$wsMap.Rows.Item(95).RowHeight = 18

# ---------------------------------------------------------------------------
# 2. "Must have properties" sheet - the same six fields, reordered
#    (image.organ_or_tissue first), appended as rows 24-29. Columns:
#    A=Full Name, B=Node Name, C=Property Name.
# ---------------------------------------------------------------------------

$wsMust.Range("A24").Value = "image.organ_or_tissue"
$wsMust.Range("B24").Value = "image"
$wsMust.Range("C24").Value = "organ_or_tissue"

$wsMust.Range("A25").Value = "MultiplexMicroscopy.tissue_fixative"
$wsMust.Range("B25").Value = "MultiplexMicroscopy"
$wsMust.Range("C25").Value = "tissue_fixative"

$wsMust.Range("A26").Value = "MultiplexMicroscopy.imaging_assay_type"
$wsMust.Range("B26").Value = "MultiplexMicroscopy"
$wsMust.Range("C26").Value = "imaging_assay_type"

$wsMust.Range("A27").Value = "proteomic.analytical_fractions"
$wsMust.Range("B27").Value = "proteomic"
$wsMust.Range("C27").Value = "analytical_fractions"

$wsMust.Range("A28").Value = "proteomic.instrument_make"
$wsMust.Range("B28").Value = "proteomic"
$wsMust.Range("C28").Value = "instrument_make"

$wsMust.Range("A29").Value = "proteomic.proteomic_design_description"
$wsMust.Range("B29").Value = "proteomic"
$wsMust.Range("C29").Value = "proteomic_design_description"

$wsMust.Rows.Item(24).RowHeight = 18
$wsMust.Rows.Item(25).RowHeight = 19
$wsMust.Rows.Item(26).RowHeight = 19
$wsMust.Rows.Item(27).RowHeight = 18
$wsMust.Rows.Item(28).RowHeight = 18
$wsMust.Rows.Item(29).RowHeight = 18

# ---------------------------------------------------------------------------
# 3. Distinctive cell formatting for the two groups of new "Full Name" cells:
#      - MultiplexMicroscopy.* rows: 14pt red Courier New
#      - proteomic.* / image.organ_or_tissue rows: 12pt dark-blue Arial
#        Unicode MS
#    The first cell of each group gets the format built up property by
#    property; every other cell in the group reuses that exact style via
#    copy/paste-special so the workbook doesn't accumulate a distinct style
#    per cell.
# ---------------------------------------------------------------------------

# --- Style A: red Courier New (MultiplexMicroscopy.* rows) ---
$redCourierSource = $wsMap.Range("D90")
$f = $redCourierSource.Font
$f.Name = "Courier New"
$f.Size = 14
$f.Color = 255        # BGR(0,0,255) -> RGB FF0000
$f.Family = 1

$redCourierTargets = @(
    $wsMap.Range("D91"),
    $wsMust.Range("A25"),
    $wsMust.Range("A26")
)
foreach ($target in $redCourierTargets) {
    $redCourierSource.Copy()
    $target.PasteSpecial(-4122)   # xlPasteFormats
}

# --- Style B: dark-blue Arial Unicode MS (proteomic.* / image rows) ---
$navyArialSource = $wsMap.Range("D92")
$f2 = $navyArialSource.Font
$f2.Name = "Arial Unicode MS"
$f2.Size = 12
$f2.Color = 5057303    # BGR(0x4D,0x2B,0x17) -> RGB 172B4D
$f2.Family = 2

$navyArialTargets = @(
    $wsMap.Range("D93"),
    $wsMap.Range("D94"),
    $wsMap.Range("D95"),
    $wsMust.Range("A24"),
    $wsMust.Range("A27"),
    $wsMust.Range("A28"),
    $wsMust.Range("A29")
)
foreach ($target in $navyArialTargets) {
    $navyArialSource.Copy()
    $target.PasteSpecial(-4122)   # xlPasteFormats
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. "Change Log" sheet - new entry for this release.
# ---------------------------------------------------------------------------

$wsLog.Range("A6").Copy()
$wsLog.Range("A7").PasteSpecial(-4122)
$wsLog.Range("A7").Value = 45637   # 2024-12-11
$wsLog.Range("B7").Value = "Bruce Wang"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Restore the active selections on each sheet to match where the author
#    left off editing.
# ---------------------------------------------------------------------------

$wsMap.Activate()
$wsMap.Range("B71").Select()

$wsMust.Activate()
$wsMust.Range("B24:C24").Select()

$wsLog.Activate()
$wsLog.Range("B17").Select()

$wsMap.Activate()
